$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to remain a text cell (matches the source data,
    # which stores prices like "23.700.30" / "0.07900" as literal strings)
    # instead of being auto-coerced to a Number by COM, then restore the
    # default "Normal" style so no stray formatting is introduced.
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "23.676.10"
$ws.Range("E2").Value = "  +0.71%  "

Set-TextValue "D3" "1.651.43"
$ws.Range("E3").Value = "  +0.60%  "

Set-TextValue "D4" "1.003"
$ws.Range("E4").Value = "  +0.34%  "

Set-TextValue "D5" "1.002"
$ws.Range("E5").Value = "  +0.21%  "

Set-TextValue "D6" "302.96"
$ws.Range("E6").Value = "  -0.54%  "

Set-TextValue "D7" "0.3801"
$ws.Range("E7").Value = "  +0.33%  "

Set-TextValue "D8" "0.3606"
$ws.Range("E8").Value = "  -0.55%  "

Set-TextValue "D9" "51.17"
$ws.Range("E9").Value = "  -1.06%  "

Set-TextValue "D10" "0.08186"
$ws.Range("E10").Value = "  +0.09%  "

Set-TextValue "D11" "1.228"
$ws.Range("E11").Value = "  -0.73%  "

Set-TextValue "D12" "1.003"
$ws.Range("E12").Value = "  +0.36%  "

Set-TextValue "D13" "22.50"
$ws.Range("E13").Value = "  -0.16%  "

Set-TextValue "D14" "6.483"
$ws.Range("E14").Value = "  +0.16%  "

Set-TextValue "D15" "7.392"
$ws.Range("E15").Value = "  +0.17%  "

Set-TextValue "D16" "0.00001229"
$ws.Range("E16").Value = "  -1.04%  "

Set-TextValue "D17" "1.651.13"
$ws.Range("E17").Value = "  +0.98%  "

Set-TextValue "D18" "97.35"
$ws.Range("E18").Value = "  +1.63%  "

Set-TextValue "D19" "0.07004"
$ws.Range("E19").Value = "  +1.06%  "

Set-TextValue "D20" "6.819"
$ws.Range("E20").Value = "  +3.42%  "

Set-TextValue "D21" "17.59"
$ws.Range("E21").Value = "  +0.33%  "

Set-TextValue "D22" "1.002"
$ws.Range("E22").Value = "  +0.27%  "

Set-TextValue "D23" "12.76"
$ws.Range("E23").Value = "  +1.98%  "

Set-TextValue "D24" "23.680.40"
$ws.Range("E24").Value = "  +0.74%  "

Set-TextValue "D25" "2.516"
$ws.Range("E25").Value = "  +0.44%  "

Set-TextValue "D26" "3.027"
$ws.Range("E26").Value = "  -1.23%  "

Set-TextValue "D27" "21.19"
$ws.Range("E27").Value = "  +0.12%  "

Set-TextValue "D28" "152.84"
$ws.Range("E28").Value = "  +0.72%  "

Set-TextValue "D29" "5.219"
$ws.Range("E29").Value = "  -0.48%  "

Set-TextValue "D30" "133.87"
$ws.Range("E30").Value = "  +0.47%  "

Set-TextValue "D31" "1.840.74"
$ws.Range("E31").Value = "  +1.17%  "

Set-TextValue "D32" "6.951"
$ws.Range("E32").Value = "  +4.66%  "

Set-TextValue "D33" "2.205"
$ws.Range("E33").Value = "  +2.56%  "

Set-TextValue "D34" "12.04"
$ws.Range("E34").Value = "  +4.87%  "

Set-TextValue "D35" "1.058"
$ws.Range("E35").Value = "  -1.26%  "

Set-TextValue "D36" "0.02796"
$ws.Range("E36").Value = "  +1.07%  "

Set-TextValue "D37" "0.2513"
$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D38" "0.08773"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D39" "6.064"
$ws.Range("E39").Value = "  +0.74%  "

Set-TextValue "D40" "0.07022"
$ws.Range("E40").Value = "  -1.11%  "

Set-TextValue "D41" "13.04"
$ws.Range("E41").Value = "  +6.92%  "

Set-TextValue "D42" "0.6995"
$ws.Range("E42").Value = "  -0.85%  "

Set-TextValue "D43" "1.336"
$ws.Range("E43").Value = "  -0.91%  "

Set-TextValue "D44" "15.92"
$ws.Range("E44").Value = "  +1.10%  "

Set-TextValue "D45" "0.6505"
$ws.Range("E45").Value = "  -0.42%  "

$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D46" "1.002"
$ws.Range("E46").Value = "  +0.35%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D47" "2.309"
$ws.Range("E47").Value = "  +1.14%  "

Set-TextValue "D48" "3.962"
$ws.Range("E48").Value = "  -0.30%  "

Set-TextValue "D49" "0.07900"
$ws.Range("E49").Value = "  -1.02%  "

Set-TextValue "D50" "127.84"
$ws.Range("E50").Value = "  -0.69%  "

Set-TextValue "D51" "1.179"
$ws.Range("E51").Value = "  -1.17%  "
